# Apply the 2022-05-11 data refresh for Fonds de solidarite (volet 1, regional, classe effectif)
# Updates nombre_aides (column C) and montant_total (column E) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C14").Value = 110815
$ws.Range("E14").Value = 253243499

$ws.Range("C38").Value = 7250
$ws.Range("E38").Value = 58549613

$ws.Range("C51").Value = 6353
$ws.Range("E51").Value = 12182169

$ws.Range("C53").Value = 141677
$ws.Range("E53").Value = 590053180

$ws.Range("C57").Value = 3710
$ws.Range("E57").Value = 138291366

$ws.Range("C64").Value = 5195
$ws.Range("E64").Value = 20316077

$ws.Range("C66").Value = 766
$ws.Range("E66").Value = 9908094

$ws.Range("C72").Value = 6273
$ws.Range("E72").Value = 15083703

$ws.Range("C79").Value = 116589
$ws.Range("E79").Value = 447328195

$ws.Range("C82").Value = 8451
$ws.Range("E82").Value = 124850808

$ws.Range("C91").Value = 151086
$ws.Range("E91").Value = 482026734

$ws.Range("C92").Value = 408967
$ws.Range("E92").Value = 1593247574

$ws.Range("C93").Value = 209456
$ws.Range("E93").Value = 1307100586

$ws.Range("C94").Value = 94127
$ws.Range("E94").Value = 915011909

$ws.Range("C95").Value = 50692
$ws.Range("E95").Value = 929113177

$ws.Range("C96").Value = 17219
$ws.Range("E96").Value = 787375015

$ws.Range("C104").Value = 135212
$ws.Range("E104").Value = 272075024

$ws.Range("C105").Value = 8164
$ws.Range("E105").Value = 16862003

$ws.Range("C106").Value = 18333
$ws.Range("E106").Value = 41276024

$ws.Range("C108").Value = 2831
$ws.Range("E108").Value = 18485757

$ws.Range("C109").Value = 1269
$ws.Range("E109").Value = 20732303

$ws.Range("C113").Value = 8803
$ws.Range("E113").Value = 12663437

$ws.Range("C115").Value = 11685
$ws.Range("E115").Value = 32906669

$ws.Range("C116").Value = 4549
$ws.Range("E116").Value = 20417858

$ws.Range("C142").Value = 168967
$ws.Range("E142").Value = 681734390

$ws.Range("C147").Value = 408
$ws.Range("E147").Value = 28868306

$ws.Range("C154").Value = 201565
$ws.Range("E154").Value = 786746320

$ws.Range("C161").Value = 36
$ws.Range("E161").Value = 3796789

$ws.Range("C163").Value = 70983
$ws.Range("E163").Value = 131766598

$ws.Range("C172").Value = 22699
$ws.Range("E172").Value = 44669168

$ws.Range("C173").Value = 96857
$ws.Range("E173").Value = 327922035

$ws.Range("C174").Value = 226079
$ws.Range("E174").Value = 900505380

$ws.Range("C175").Value = 80777
$ws.Range("E175").Value = 485996075

$wb.Save()
